$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of receipt data (row 4). Column A holds a date-formatted string
# that must stay as literal text (matching the existing rows), so force
# the cell to Text format before assigning the value, then restore the
# default "Normal" style so no stray formatting is left behind.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2016-09-29"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = "FamilyMart 麻布十番一丁目店"
$ws.Range("C4").Value = "茎わかめ、ハムカツサンド、中華そば"
$ws.Range("D4").Value = 512
